# Generate Report for Handoff
# Updates the localization-status report: the "Handed back: in sync with en-US"
# status becomes "Ready for handoff" (the handoff-to-loc step just ran), and the
# associated timestamps advance a couple of minutes. Touching those cells makes
# the (now shorter) Status column narrower, so we also shrink the Status columns
# on all three sheets to match their new best-fit width.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-11-09 01:11:51"
$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-11-09 01:11:38"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.25

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-11-09 01:11:51"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.25
